$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy number formats from the row above (row 12) to keep styles consistent
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$ws.Range("A13").Value = 44079
$ws.Range("B13").Value = 0.83333333333333337
$ws.Range("C13").Value = 1.03125
$ws.Range("E13").Value = "xml toimii unityssa"
$ws.Range("E14").Value = "database uudelleenmaarittely + unity"

$ws.Range("D16").Select()
